$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.067.40'
$ws.Range('E2').Value = '  +1.74%  '
$ws.Range('D3').Value = '1.856.21'
$ws.Range('E3').Value = '  +3.21%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.32'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.623'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.83%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.51'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +8.60%  '
$ws.Range('E9').Value = '  +3.47%  '
$ws.Range('E10').Value = '  +2.64%  '
$ws.Range('D12').Value = '2.126.46'
$ws.Range('E12').Value = '  +3.25%  '
$ws.Range('D13').Value = '1.871.25'
$ws.Range('E13').Value = '  +3.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.42'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.677'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.69'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +3.52%  '
$ws.Range('D17').Value = '35.058.05'
$ws.Range('E17').Value = '  +2.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.26'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.96%  '
$ws.Range('D19').Value = '0.0₃0795'
$ws.Range('E19').Value = '  +2.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '240.84'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.13'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.32%  '
$ws.Range('E22').Value = '  +1.43%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('E24').Value = '  +1.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.10'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.99%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.93'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +29.96%  '
$ws.Range('E27').Value = '  +3.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.69'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.22%  '
$ws.Range('E29').Value = '  +2.43%  '
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('E31').Value = '  +3.31%  '
$ws.Range('E32').Value = '  +0.46%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.00'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +3.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.02'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +13.38%  '
$ws.Range('E35').Value = '  +22.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.786'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +13.86%  '
$ws.Range('E37').Value = '  +6.13%  '
$ws.Range('E38').Value = '  +13.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '91.78'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.58%  '
$ws.Range('E40').Value = '  +7.05%  '
$ws.Range('D41').Value = '1.352.01'
$ws.Range('E41').Value = '  +2.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.85'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +5.23%  '
$ws.Range('E43').Value = '  +6.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.86'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +56.73%  '
$ws.Range('E45').Value = '  +0.56%  '
$ws.Range('E46').Value = '  +2.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0540'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +5.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.43'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +4.90%  '
$ws.Range('D49').Value = '2.041.04'
$ws.Range('E50').Value = '  +3.20%  '
$ws.Range('E51').Value = '  +18.38%  '
